$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 21739634
$ws.Range("I92").Value = 22727712
$ws.Range("K92").Value = 22727712
$ws.Range("M92").Value = -22726464

$ws.Range("H137").Value = 2008.1428
$ws.Range("I137").Value = 858.3333
$ws.Range("K137").Value = 2574.9999
$ws.Range("M137").Value = -24.9998999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1948.7037
$ws.Range("I2").Value = 1405.1
$ws.Range("J2").Value = 3501.8572
$ws.Range("K2").Value = 1405.1
$ws.Range("L2").Value = 3501.8572
$ws.Range("M2").Value = -1292.1
$ws.Range("N2").Value = -3727.8572

$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

$ws.Range("H45").Value = 4136.3335
$ws.Range("I45").Value = 2954.75
$ws.Range("K45").Value = 2954.75
$ws.Range("M45").Value = -2577.75

$ws.Range("H61").Value = 3109.0334
$ws.Range("I61").Value = 2955.1924
$ws.Range("K61").Value = 2955.1924
$ws.Range("M61").Value = -2743.1924

$ws.Range("H74").Value = 2501.45
$ws.Range("I74").Value = 1536
$ws.Range("J74").Value = 5397.8
$ws.Range("K74").Value = 1536
$ws.Range("L74").Value = 5397.8
$ws.Range("M74").Value = -662
$ws.Range("N74").Value = -7145.8

$ws.Range("H77").Value = 2501.45
$ws.Range("I77").Value = 1536
$ws.Range("J77").Value = 5397.8
$ws.Range("K77").Value = 7680
$ws.Range("L77").Value = 26989
$ws.Range("M77").Value = -3312
$ws.Range("N77").Value = -35725

$ws.Range("H116").Value = 1948.7037
$ws.Range("I116").Value = 1405.1
$ws.Range("J116").Value = 3501.8572
$ws.Range("K116").Value = 1405.1
$ws.Range("L116").Value = 3501.8572
$ws.Range("M116").Value = 888.9000000000001
$ws.Range("N116").Value = -8089.8572

$ws.Range("H136").Value = 3109.0334
$ws.Range("I136").Value = 2955.1924
$ws.Range("K136").Value = 8865.5772
$ws.Range("M136").Value = -6315.5772

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1948.7037
$ws.Range("I3").Value = 1405.1
$ws.Range("J3").Value = 3501.8572
$ws.Range("K3").Value = 1405.1
$ws.Range("L3").Value = 3501.8572
$ws.Range("M3").Value = -1291.1
$ws.Range("N3").Value = -3729.8572

$ws.Range("H82").Value = 39234.855
$ws.Range("I82").Value = 7084.6665
$ws.Range("J82").Value = 63347.5
$ws.Range("K82").Value = 7084.6665
$ws.Range("L82").Value = 63347.5
$ws.Range("M82").Value = -6701.6665
$ws.Range("N82").Value = -64113.5

$ws.Range("H85").Value = 39234.855
$ws.Range("I85").Value = 7084.6665
$ws.Range("J85").Value = 63347.5
$ws.Range("K85").Value = 7084.6665
$ws.Range("L85").Value = 63347.5
$ws.Range("M85").Value = -5758.6665
$ws.Range("N85").Value = -65999.5

$ws.Range("H107").Value = 903.8261
$ws.Range("I107").Value = 777.75
$ws.Range("K107").Value = 777.75
$ws.Range("M107").Value = 1142.25

$ws.Range("H134").Value = 22491904
$ws.Range("I134").Value = 5497764
$ws.Range("K134").Value = 16493292
$ws.Range("M134").Value = -16490757

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3571.8845
$ws.Range("I31").Value = 2212.625
$ws.Range("J31").Value = 5746.7
$ws.Range("K31").Value = 2212.625
$ws.Range("L31").Value = 5746.7
$ws.Range("M31").Value = -1917.625
$ws.Range("N31").Value = -6336.7

$ws.Range("H34").Value = 3571.8845
$ws.Range("I34").Value = 2212.625
$ws.Range("J34").Value = 5746.7
$ws.Range("K34").Value = 2212.625
$ws.Range("L34").Value = 5746.7
$ws.Range("M34").Value = -2010.625
$ws.Range("N34").Value = -6150.7

$ws.Range("H110").Value = 89995
$ws.Range("J110").Value = 89995
$ws.Range("L110").Value = 89995
$ws.Range("N110").Value = -98175

$ws.Range("H132").Value = 1870
$ws.Range("I132").Value = 1801.1538
$ws.Range("K132").Value = 5403.4614
$ws.Range("M132").Value = -2873.4614

$ws.Range("H134").Value = 3687.5715
$ws.Range("J134").Value = 4888.1665
$ws.Range("L134").Value = 14664.4995
$ws.Range("N134").Value = -19734.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 5500
$ws.Range("I58").Value = 5500
$ws.Range("K58").Value = 16500
$ws.Range("M58").Value = -16372

$ws.Range("H131").Value = 1664.9048
$ws.Range("I131").Value = 1332.5714
$ws.Range("J131").Value = 1831.0714
$ws.Range("K131").Value = 3997.7142
$ws.Range("L131").Value = 5493.2142
$ws.Range("M131").Value = 1042.2858
$ws.Range("N131").Value = -15573.2142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4460.12
$ws.Range("J70").Value = 4500.174
$ws.Range("L70").Value = 4500.174
$ws.Range("N70").Value = -5040.174

$ws.Range("H73").Value = 4460.12
$ws.Range("J73").Value = 4500.174
$ws.Range("L73").Value = 4500.174
$ws.Range("N73").Value = -6372.174

$ws.Range("H102").Value = 3149.9
$ws.Range("I102").Value = 2800.4
$ws.Range("J102").Value = 3499.4
$ws.Range("K102").Value = 2800.4
$ws.Range("L102").Value = 3499.4
$ws.Range("M102").Value = -1178.4
$ws.Range("N102").Value = -6743.4

$ws.Range("H132").Value = 2902.1086
$ws.Range("I132").Value = 1938.0938
$ws.Range("K132").Value = 5814.2814
$ws.Range("M132").Value = -3284.2814

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5614.8335
$ws.Range("I7").Value = 3939.8572
$ws.Range("K7").Value = 3939.8572
$ws.Range("M7").Value = -3827.8572

$ws.Range("H16").Value = 1542.0526
$ws.Range("I16").Value = 1370.5883
$ws.Range("K16").Value = 1370.5883
$ws.Range("M16").Value = -1200.5883

$ws.Range("H22").Value = 1982.3889
$ws.Range("I22").Value = 1812.4166
$ws.Range("J22").Value = 2322.3333
$ws.Range("K22").Value = 1812.4166
$ws.Range("L22").Value = 2322.3333
$ws.Range("M22").Value = -1517.4166
$ws.Range("N22").Value = -2912.3333

$ws.Range("H27").Value = 1982.3889
$ws.Range("I27").Value = 1812.4166
$ws.Range("J27").Value = 2322.3333
$ws.Range("K27").Value = 1812.4166
$ws.Range("L27").Value = 2322.3333
$ws.Range("M27").Value = -1705.4166
$ws.Range("N27").Value = -2536.3333

$ws.Range("H46").Value = 3880.739
$ws.Range("I46").Value = 1499
$ws.Range("J46").Value = 4238
$ws.Range("K46").Value = 1499
$ws.Range("L46").Value = 4238
$ws.Range("M46").Value = -1311
$ws.Range("N46").Value = -4614

$ws.Range("H126").Value = 5614.8335
$ws.Range("I126").Value = 3939.8572
$ws.Range("K126").Value = 11819.5716
$ws.Range("M126").Value = -9349.571599999999

$ws.Range("H132").Value = 2589.739
$ws.Range("I132").Value = 2272.158
$ws.Range("J132").Value = 4098.25
$ws.Range("K132").Value = 6816.474
$ws.Range("L132").Value = 12294.75
$ws.Range("M132").Value = -4286.474
$ws.Range("N132").Value = -17354.75

$ws.Range("H136").Value = 3421.325
$ws.Range("I136").Value = 2944.1072
$ws.Range("K136").Value = 8832.321599999999
$ws.Range("M136").Value = -6282.321599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3290.2974
$ws.Range("I122").Value = 2822.862
$ws.Range("J122").Value = 4984.75
$ws.Range("K122").Value = 8468.585999999999
$ws.Range("L122").Value = 14954.25
$ws.Range("M122").Value = -6018.585999999999
$ws.Range("N122").Value = -19854.25

$ws.Range("H132").Value = 2216.2856
$ws.Range("I132").Value = 1705.2142
$ws.Range("K132").Value = 5115.642599999999
$ws.Range("M132").Value = -2585.642599999999

$ws.Range("H136").Value = 2774.204
$ws.Range("I136").Value = 1885.1818
$ws.Range("J136").Value = 4607.8125
$ws.Range("K136").Value = 5655.5454
$ws.Range("L136").Value = 13823.4375
$ws.Range("M136").Value = -3105.5454
$ws.Range("N136").Value = -18923.4375
